$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ------------------------------------------------------------------
# 1. Duplicate an existing text box ("ZoneTexte 1") so the new shape
#    inherits proper run formatting (lang="fr-FR" dirty="0"), then
#    reposition / rename / retext it into "ZoneTexte 2".
# ------------------------------------------------------------------
$srcIndex = 0
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "ZoneTexte 1") {
        $srcIndex = $i
    }
}
$src = $s.Shapes.Item($srcIndex)
$dupRange = $src.Duplicate()
$newShape = $dupRange.Item(1)
$newShape.Name = "ZoneTexte 2"

$eps = 0.00001
$newShape.Left = (9105900 / 12700.0) + $eps
$newShape.Top = (1895475 / 12700.0) + $eps
$newShape.Width = (2552700 / 12700.0) + $eps
$newShape.Height = (2862322 / 12700.0) + $eps

$tr = $newShape.TextFrame.TextRange
$tr.Text = "1- Moteur Att. 1"
$tr.InsertAfter("`r2-+5V photo-transistor") | Out-Null

# Split the second paragraph's text into two runs, matching the
# original "2-+5V " / "photo-transistor" split.
$total = $tr.Text.Length
$para2Len = "2-+5V photo-transistor".Length
$startOfPhoto = $total - $para2Len + "2-+5V ".Length + 1
$sub = $tr.Characters($startOfPhoto, "photo-transistor".Length)
$sub.Text = "photo-transistor"

$tr.InsertAfter("`rAvec 100K ohms R et mesure tension //") | Out-Null
$tr.InsertAfter("`r3- Moteur Att. 2") | Out-Null
$tr.InsertAfter("`r4-+5V phototransistor avec 390 ohms R ") | Out-Null
$tr.InsertAfter("`r5- Moteur Att. 3") | Out-Null
$tr.InsertAfter("`r6- GND transistor ") | Out-Null
$tr.InsertAfter("`r7-Moteur Att. 4 ") | Out-Null

# ------------------------------------------------------------------
# 2. Select every shape on the slide (the 18 original ones plus the
#    new "ZoneTexte 2") and group them into "Groupe 20".
# ------------------------------------------------------------------
$names = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $names += $s.Shapes.Item($i).Name
}
$range = $s.Shapes.Range($names)
$grp = $range.Group()
$grp.Name = "Groupe 20"
